$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the "datetimeFigureOut" date placeholder (25.04.2017 ->
#    28.04.2017) on every slide layout and on the slide master.
# ---------------------------------------------------------------------
$newDate = "28.04.2017"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    foreach ($shp in $layout.Shapes) {
        if ($shp.Name -like "Datumsplatzhalter*" -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

foreach ($shp in $master.Shapes) {
    if ($shp.Name -like "Datumsplatzhalter*" -and $shp.HasTextFrame) {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# ---------------------------------------------------------------------
# 2) Nudge the picture on slide 1 slightly to the left.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$pic = $slide1.Shapes.Item("Grafik 5")
$pic.Left = 599.2

# ---------------------------------------------------------------------
# 3) Slide 2: terminate the first bullet with a period.
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$body = $slide2.Shapes.Item("Inhaltsplatzhalter 2")
$tr = $body.TextFrame.TextRange
$firstLine = $tr.Characters(1, 39)
[void]$firstLine.InsertAfter(". ")
